$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F: validation-rule text for "Ias39Loss.LosCod" field changed
#     from "numeric 1-99" rule to "code per dropdown menu" rule.
$ws.Range("F21").Value = "1.限輸入代碼，檢核條件：依選單/V(H)2.Ias39Loss.LosCod"
$ws.Range("F31").Value = "1.自動顯示原值，限輸入代碼，檢核條件：依選單/V(H)2.Ias39Loss.LosCod"

# --- Column M (rows 12-31): reference-document text bumped from
#     requirement-spec version V1.9 to V1.93.
$refDoc = "製作依據之需求規格書與版本：PJ201800012_URS_7介接外部系統_V1.93.docx"
for ($r = 12; $r -le 31; $r++) {
    $ws.Range("M$r").Value = $refDoc
}

# --- Column Q (rows 12-31): test-case creation date moved from
#     2022/1/17 (serial 44578) to 2022/2/16 (serial 44608).
for ($r = 12; $r -le 31; $r++) {
    $ws.Range("Q$r").Value = 44608
}
